# Add Denmark, Norway and Sweden test-data sheets (copies of the Belgium
# sheet, trimmed down to a single data row) and update the market name /
# panel cells on each.

$wb = $excel.ActiveWorkbook
$belgium = $wb.Worksheets.Item("Belgium")

# --- Denmark ----------------------------------------------------------
$belgium.Copy($null, $belgium)
$denmark = $wb.Worksheets.Item($belgium.Index + 1)
$denmark.Name = "Denmark"
$denmark.Range("A9:A10").EntireRow.Delete() | Out-Null
$denmark.Range("B2").Value = "Denmark Market"
$denmark.Range("B4").ClearContents() | Out-Null

# --- Sweden (filled in second, so its new string lands before Norway's)
$denmark.Copy($null, $denmark)
$sweden = $wb.Worksheets.Item($denmark.Index + 1)
$sweden.Name = "Sweden"
$sweden.Range("B2").Value = "Sweden Market"
$sweden.Range("A8").Value = "FC602S"

# --- Norway (created third, then dragged to sit before Sweden) --------
$sweden.Copy($null, $sweden)
$norway = $wb.Worksheets.Item($sweden.Index + 1)
$norway.Name = "Norway"
$norway.Range("B2").Value = "Norway Market"
$norway.Range("A8").Value = "FC604S"

$norway.Move($sweden)

# Sheet objects become stale once the tab order is reshuffled by Move, so
# re-resolve everything we still need to touch by name.
$belgium = $wb.Worksheets.Item("Belgium")
$denmark = $wb.Worksheets.Item("Denmark")
$sweden = $wb.Worksheets.Item("Sweden")
$norway = $wb.Worksheets.Item("Norway")

# --- Selections / active sheet -----------------------------------------
$belgium.Activate()
$belgium.Range("A1:XFD1048576").Select() | Out-Null

$denmark.Activate()
$denmark.Range("A1:XFD1048576").Select() | Out-Null

$sweden.Activate()
$sweden.Range("B4").Select() | Out-Null

$norway.Activate()
$norway.Range("B5").Select() | Out-Null
